$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. RESPONSE TIME (column G) - refresh the numeric-looking values.
#    These are stored as text (shared-string) cells in the original
#    sheet, so force a Text number format before writing, then reset
#    the style back to Normal so no stray style index is left behind.
# ------------------------------------------------------------------
$gValues = @{
    2  = "3784"
    3  = "8431"
    4  = "7354"
    5  = "8288"
    6  = "2866"
    7  = "2979"
    8  = "2926"
    9  = "3197"
    10 = "2919"
    11 = "3021"
    12 = "3127"
    13 = "3071"
    14 = "3085"
    15 = "3003"
    16 = "2967"
    17 = "3816"
    18 = "2275"
    19 = "2059"
    20 = "1994"
    21 = "1917"
}

$gRange = $ws.Range("G2:G21")
$gRange.NumberFormat = "@"
foreach ($r in 2..21) {
    $ws.Range("G$r").Value = $gValues[$r]
}
$gRange.Style = "Normal"

# ------------------------------------------------------------------
# 2. ENVIRONMENT (column O) - every data row now reports "klif"
#    instead of being blank.
# ------------------------------------------------------------------
foreach ($r in 2..21) {
    $ws.Range("O$r").Value = "klif"
}

# ------------------------------------------------------------------
# 3. New VERSION column (Q) - header + "v1" for every data row.
#    Copy the header formatting from the existing O1 header cell so
#    the new header matches the look of the others.
# ------------------------------------------------------------------
$ws.Range("O1").Copy() | Out-Null
$ws.Range("Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("Q1").Value = "VERSION"

foreach ($r in 2..21) {
    $ws.Range("Q$r").Value = "v1"
}

# Column Q should carry the same width as the other data columns.
$ws.Columns.Item(17).ColumnWidth = 34.17

Write-Output "done"
